$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 223, shifting rows 223:273 down to 224:274
$ws.Range("A223").EntireRow.Insert()

# Populate the new row 223 with the new data point
$ws.Range("A223").Value = 10
$ws.Range("B223").Value = "Vega Modelo de Temuco"
$ws.Range("C223").Value = "La Araucanía"
$ws.Range("D223").Value = 44855
$ws.Range("E223").Value = 9
$ws.Range("F223").Value = 100112052
$ws.Range("G223").Value = "Albahaca"
$ws.Range("H223").Value = "Sin especificar"
$ws.Range("I223").Value = "Primera"
$ws.Range("J223").Value = 20
$ws.Range("K223").Value = 8000
$ws.Range("L223").Value = 8000
$ws.Range("M223").Value = 8000
$ws.Range("N223").Value = "$/paquete"
$ws.Range("O223").Value = "Región de Arica y Parinacota"
$ws.Range("P223").Value = 8000
$ws.Range("Q223").Value = 1
$ws.Range("R223").Value = "Hortaliza"
